# Auto-generated edit script: applies numeric cell updates to the
# "Moogle_Profits" leve-crafting profit tables across multiple sheets.
# Each sheet has columns H..N holding market-price/profit figures that
# were refreshed by the scheduled data-pull runner.
$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2262
$ws.Range("J17").Value = 2262
$ws.Range("L17").Value = 6786
$ws.Range("N17").Value = -7122
$ws.Range("H87").Value = 58077.668
$ws.Range("J87").Value = 59962.5
$ws.Range("L87").Value = 59962.5
$ws.Range("N87").Value = -62458.5
$ws.Range("H90").Value = 58077.668
$ws.Range("J90").Value = 59962.5
$ws.Range("L90").Value = 179887.5
$ws.Range("N90").Value = -192367.5
$ws.Range("H113").Value = 4597.9473
$ws.Range("I113").Value = 3828.4546
$ws.Range("J113").Value = 5656
$ws.Range("K113").Value = 3828.4546
$ws.Range("L113").Value = 5656
$ws.Range("M113").Value = -574.4546
$ws.Range("N113").Value = -12164
$ws.Range("H132").Value = 4076.0667
$ws.Range("I132").Value = 4241.231
$ws.Range("K132").Value = 12723.693
$ws.Range("M132").Value = -10193.693
$ws.Range("H137").Value = 1446.6086
$ws.Range("I137").Value = 1359.3143
$ws.Range("J137").Value = 1724.3636
$ws.Range("K137").Value = 4077.9429
$ws.Range("L137").Value = 5173.0908
$ws.Range("M137").Value = -1527.9429
$ws.Range("N137").Value = -10273.0908
$ws.Range("H138").Value = 2966.8298
$ws.Range("I138").Value = 2671.5
$ws.Range("J138").Value = 3302.432
$ws.Range("K138").Value = 8014.5
$ws.Range("L138").Value = 9907.295999999998
$ws.Range("M138").Value = -2874.5
$ws.Range("N138").Value = -20187.296

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4903.4287
$ws.Range("I32").Value = 2300.08
$ws.Range("K32").Value = 2300.08
$ws.Range("M32").Value = -2013.08
$ws.Range("H62").Value = 45083
$ws.Range("J62").Value = 45083
$ws.Range("L62").Value = 45083
$ws.Range("N62").Value = -46331
$ws.Range("H65").Value = 45083
$ws.Range("J65").Value = 45083
$ws.Range("L65").Value = 135249
$ws.Range("N65").Value = -141489
$ws.Range("H74").Value = 4637.724
$ws.Range("I74").Value = 1791.5883
$ws.Range("J74").Value = 8669.75
$ws.Range("K74").Value = 1791.5883
$ws.Range("L74").Value = 8669.75
$ws.Range("M74").Value = -917.5882999999999
$ws.Range("N74").Value = -10417.75
$ws.Range("H77").Value = 4637.724
$ws.Range("I77").Value = 1791.5883
$ws.Range("J77").Value = 8669.75
$ws.Range("K77").Value = 8957.941499999999
$ws.Range("L77").Value = 43348.75
$ws.Range("M77").Value = -4589.941499999999
$ws.Range("N77").Value = -52084.75
$ws.Range("H92").Value = 62000
$ws.Range("J92").Value = 62000
$ws.Range("L92").Value = 62000
$ws.Range("N92").Value = -66992
$ws.Range("H97").Value = 2402.0908
$ws.Range("I97").Value = 2452.3
$ws.Range("J97").Value = 1900
$ws.Range("K97").Value = 2452.3
$ws.Range("L97").Value = 1900
$ws.Range("M97").Value = -1956.3
$ws.Range("N97").Value = -2892
$ws.Range("H102").Value = 9998.333000000001
$ws.Range("H132").Value = 4575.0347
$ws.Range("I132").Value = 2057.7368
$ws.Range("K132").Value = 6173.2104
$ws.Range("M132").Value = -3643.2104

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 438.01923
$ws.Range("I7").Value = 464.6129
$ws.Range("K7").Value = 464.6129
$ws.Range("M7").Value = -351.6129
$ws.Range("H31").Value = 5057.0435
$ws.Range("I31").Value = 2357.9656
$ws.Range("J31").Value = 9661.352999999999
$ws.Range("K31").Value = 2357.9656
$ws.Range("L31").Value = 9661.352999999999
$ws.Range("M31").Value = -2062.9656
$ws.Range("N31").Value = -10251.353
$ws.Range("H34").Value = 5057.0435
$ws.Range("I34").Value = 2357.9656
$ws.Range("J34").Value = 9661.352999999999
$ws.Range("K34").Value = 2357.9656
$ws.Range("L34").Value = 9661.352999999999
$ws.Range("M34").Value = -2155.9656
$ws.Range("N34").Value = -10065.353
$ws.Range("H45").Value = 64000
$ws.Range("J45").Value = 64000
$ws.Range("L45").Value = 64000
$ws.Range("N45").Value = -65186
$ws.Range("H62").Value = 3528.0588
$ws.Range("I62").Value = 3076
$ws.Range("K62").Value = 3076
$ws.Range("M62").Value = -2452
$ws.Range("H65").Value = 3528.0588
$ws.Range("I65").Value = 3076
$ws.Range("K65").Value = 15380
$ws.Range("M65").Value = -12260
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 3950.875
$ws.Range("I105").Value = 3372.4285
$ws.Range("J105").Value = 8000
$ws.Range("K105").Value = 3372.4285
$ws.Range("L105").Value = 8000
$ws.Range("M105").Value = -1625.4285
$ws.Range("N105").Value = -11494
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 72.07143000000001
$ws.Range("I6").Value = 50.727272
$ws.Range("K6").Value = 152.181816
$ws.Range("M6").Value = -39.181816
$ws.Range("H40").Value = 638.9474
$ws.Range("I40").Value = 687.17645
$ws.Range("K40").Value = 2748.7058
$ws.Range("M40").Value = -2679.7058
$ws.Range("H52").Value = 765.5
$ws.Range("J52").Value = 765.5
$ws.Range("L52").Value = 2296.5
$ws.Range("N52").Value = -2828.5
$ws.Range("H55").Value = 3437.25
$ws.Range("I55").Value = 1266.3334
$ws.Range("J55").Value = 9950
$ws.Range("K55").Value = 3799.0002
$ws.Range("L55").Value = 29850
$ws.Range("M55").Value = -3622.0002
$ws.Range("N55").Value = -30204
$ws.Range("H62").Value = 5644.909
$ws.Range("I62").Value = 5813.4
$ws.Range("J62").Value = 5504.5
$ws.Range("K62").Value = 17440.2
$ws.Range("L62").Value = 16513.5
$ws.Range("M62").Value = -16754.2
$ws.Range("N62").Value = -17885.5
$ws.Range("H65").Value = 5644.909
$ws.Range("I65").Value = 5813.4
$ws.Range("J65").Value = 5504.5
$ws.Range("K65").Value = 52320.6
$ws.Range("L65").Value = 49540.5
$ws.Range("M65").Value = -48888.6
$ws.Range("N65").Value = -56404.5
$ws.Range("H95").Value = 9999
$ws.Range("I95").Value = 9999
$ws.Range("K95").Value = 29997
$ws.Range("M95").Value = -27938
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 1000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -754
$ws.Range("H128").Value = 511099
$ws.Range("I128").Value = 511099
$ws.Range("K128").Value = 1533297
$ws.Range("M128").Value = -1528317
$ws.Range("H139").Value = 2388.077
$ws.Range("I139").Value = 1362.8
$ws.Range("K139").Value = 4088.4
$ws.Range("M139").Value = 1051.6

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 797.1429000000001
$ws.Range("I97").Value = 461.33334
$ws.Range("K97").Value = 461.33334
$ws.Range("M97").Value = 34.66665999999998
$ws.Range("H113").Value = 4990.1177
$ws.Range("I113").Value = 4254.8184
$ws.Range("K113").Value = 4254.8184
$ws.Range("M113").Value = -2084.8184
$ws.Range("H123").Value = 59000
$ws.Range("J123").Value = 59000
$ws.Range("L123").Value = 59000
$ws.Range("N123").Value = -63900
$ws.Range("H132").Value = 5168.033
$ws.Range("I132").Value = 3101.8462
$ws.Range("J132").Value = 6748.0586
$ws.Range("K132").Value = 9305.5386
$ws.Range("L132").Value = 20244.1758
$ws.Range("M132").Value = -6775.5386
$ws.Range("N132").Value = -25304.1758

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 128715
$ws.Range("I7").Value = 264610.75
$ws.Range("J7").Value = 19998.4
$ws.Range("K7").Value = 264610.75
$ws.Range("L7").Value = 19998.4
$ws.Range("M7").Value = -264498.75
$ws.Range("N7").Value = -20222.4
$ws.Range("H22").Value = 1118.9231
$ws.Range("I22").Value = 512.5
$ws.Range("J22").Value = 2089.2
$ws.Range("K22").Value = 512.5
$ws.Range("L22").Value = 2089.2
$ws.Range("M22").Value = -217.5
$ws.Range("N22").Value = -2679.2
$ws.Range("H27").Value = 1118.9231
$ws.Range("I27").Value = 512.5
$ws.Range("J27").Value = 2089.2
$ws.Range("K27").Value = 512.5
$ws.Range("L27").Value = 2089.2
$ws.Range("M27").Value = -405.5
$ws.Range("N27").Value = -2303.2
$ws.Range("H88").Value = 9990
$ws.Range("I88").Value = 9990
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 9990
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -9562
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 9990
$ws.Range("I91").Value = 9990
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 9990
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -8508
$ws.Range("N91").ClearContents()
$ws.Range("H93").Value = 4453.2856
$ws.Range("I93").Value = 4714.1665
$ws.Range("K93").Value = 4714.1665
$ws.Range("M93").Value = -3466.1665
$ws.Range("H126").Value = 128715
$ws.Range("I126").Value = 264610.75
$ws.Range("J126").Value = 19998.4
$ws.Range("K126").Value = 793832.25
$ws.Range("L126").Value = 59995.2
$ws.Range("M126").Value = -791362.25
$ws.Range("N126").Value = -64935.2

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 13175
$ws.Range("J45").Value = 13175
$ws.Range("L45").Value = 13175
$ws.Range("N45").Value = -14157
$ws.Range("H96").Value = 7003.6
$ws.Range("J96").Value = 8329.5
$ws.Range("L96").Value = 8329.5
$ws.Range("N96").Value = -11075.5
$ws.Range("H126").Value = 2498.375
$ws.Range("I126").Value = 2433.6
$ws.Range("J126").Value = 2606.3333
$ws.Range("K126").Value = 7300.799999999999
$ws.Range("L126").Value = 7818.999899999999
$ws.Range("M126").Value = -4830.799999999999
$ws.Range("N126").Value = -12758.9999

